# Updating filtered feeds from workflow
# Appends one new row (row 32) to the feed table on the active sheet,
# mirroring the existing rows: a hyperlinked URL in column A, a keyword
# in column B, and a title in column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$link = "https://www.360dx.com/cancer/gestalt-primaa-partner-ai-based-testing-skin-cancer"
$keywords = "digital pathology"
$title = "Gestalt, Primaa Partner on AI-Based Testing for Skin Cancer"

$row = 32

$ws.Cells.Item($row, 1).Value = $link
$ws.Cells.Item($row, 2).Value = $keywords
$ws.Cells.Item($row, 3).Value = $title

# Turn the link cell into a real hyperlink, then restore the same
# "Hyperlink" cell style already used by the other rows in column A.
$ws.Hyperlinks.Add($ws.Cells.Item($row, 1), $link)
$ws.Cells.Item($row, 1).Style = $ws.Cells.Item($row - 1, 1).Style
